$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Electric Utilities by Census Division and State, November 2016"

$ws.Range("C4").Value = 27
$ws.Range("E4").Value = 0
$ws.Range("H4").Value = 41
$ws.Range("C5").Value = 33
$ws.Range("H5").Value = 327
$ws.Range("C6").Value = 405
$ws.Range("C7").Value = 64
$ws.Range("E7").Value = 0
$ws.Range("H7").Value = 112
$ws.Range("C8").Value = 234
$ws.Range("H8").Value = 32
$ws.Range("C9").Value = 32
$ws.Range("C10").Value = 561
$ws.Range("H10").Value = 75
$ws.Range("C11").Value = 131
$ws.Range("E11").Value = 11
$ws.Range("H11").Value = 1
$ws.Range("C12").Value = 719
$ws.Range("E12").Value = 270
$ws.Range("C13").Value = 133
$ws.Range("E13").Value = 11
$ws.Range("C14").Value = 176
$ws.Range("E14").Value = 0
$ws.Range("H14").Value = 197
$ws.Range("C15").Value = 5
$ws.Range("F15").Value = 36
$ws.Range("H15").Value = 16
$ws.Range("C16").Value = 51
$ws.Range("E16").Value = 0
$ws.Range("H16").Value = 240
$ws.Range("C17").Value = 5
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 941
$ws.Range("H17").Value = 15
$ws.Range("C18").Value = 14
$ws.Range("H18").Value = 31
$ws.Range("C19").Value = 3
$ws.Range("E19").Value = 1
$ws.Range("H19").Value = 6
$ws.Range("C20").Value = 23
$ws.Range("E20").Value = 2
$ws.Range("H20").Value = 25
$ws.Range("C21").Value = 9
$ws.Range("H21").Value = 8
$ws.Range("C22").Value = 20
$ws.Range("E22").Value = 15
$ws.Range("H22").Value = 37
$ws.Range("C23").Value = 26
$ws.Range("E23").Value = 13
$ws.Range("C24").Value = 56
$ws.Range("E24").Value = 3
$ws.Range("H24").Value = 48
$ws.Range("E25").Value = 14
$ws.Range("H25").Value = 13
$ws.Range("C26").Value = 128
$ws.Range("E26").Value = 8
$ws.Range("H26").Value = 28
$ws.Range("C27").Value = 10
$ws.Range("E27").Value = 43
$ws.Range("C28").Value = 532
$ws.Range("E28").Value = 14
$ws.Range("H28").Value = 0.44
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 9
$ws.Range("E29").Value = 0.31
$ws.Range("H29").Value = 11
$ws.Range("C30").Value = 683
$ws.Range("E30").Value = 141
$ws.Range("C31").Value = 8
$ws.Range("E31").Value = 1
$ws.Range("H31").Value = 92
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 82
$ws.Range("E32").Value = 0
$ws.Range("H32").Value = 15
$ws.Range("C33").Value = 47
$ws.Range("C34").Value = 35
$ws.Range("H34").Value = 17
$ws.Range("C35").Value = 78
$ws.Range("E35").Value = 0
$ws.Range("H35").Value = 28
$ws.Range("C36").Value = 30
$ws.Range("E36").Value = 0.07000000000000000666
$ws.Range("H36").Value = 38
$ws.Range("H37").Value = 69
$ws.Range("B38").Value = 0.28000000000000003
$ws.Range("C38").Value = 5
$ws.Range("E38").Value = 2
$ws.Range("H38").Value = 7
$ws.Range("E39").Value = 6
$ws.Range("H39").Value = 11
$ws.Range("B40").Value = 1
$ws.Range("E40").Value = 0
$ws.Range("H40").Value = 7
$ws.Range("C41").Value = 437
$ws.Range("C42").Value = 0.26
$ws.Range("H42").Value = 12
$ws.Range("H43").Value = 14
$ws.Range("E44").Value = 3
$ws.Range("H44").Value = 16
$ws.Range("C46").Value = 6
$ws.Range("H46").Value = 30
$ws.Range("E47").Value = 2
$ws.Range("H47").Value = 41
$ws.Range("C48").Value = 12
$ws.Range("H48").Value = 4
$ws.Range("C49").Value = 9
$ws.Range("C50").Value = 87
$ws.Range("E50").Value = 0
$ws.Range("H50").Value = 36
$ws.Range("C51").Value = 475
$ws.Range("E51").Value = 30
$ws.Range("H51").Value = 10
$ws.Range("B52").Value = 162
$ws.Range("C52").Value = 1551
$ws.Range("E52").Value = 62
$ws.Range("H52").Value = 4
$ws.Range("C54").Value = 56
$ws.Range("E54").Value = 7
$ws.Range("H54").Value = 116
$ws.Range("C55").Value = 4
$ws.Range("E55").Value = 7
$ws.Range("H55").Value = 42
$ws.Range("E56").Value = 140
$ws.Range("H56").Value = 27
$ws.Range("C57").Value = 32
$ws.Range("E57").Value = 3
$ws.Range("C58").Value = 12
$ws.Range("E58").Value = 4
$ws.Range("H58").Value = 6
$ws.Range("E59").Value = 0.2
$ws.Range("C60").Value = 324
$ws.Range("E60").Value = 5
$ws.Range("C61").Value = 5
$ws.Range("E61").Value = 13
$ws.Range("H61").Value = 20
$ws.Range("C62").Value = 6
$ws.Range("E62").Value = 13
$ws.Range("H62").Value = 20
$ws.Range("C63").Value = 6
$ws.Range("H63").Value = 244
$ws.Range("B64").Value = 0.22
$ws.Range("C64").Value = 4
$ws.Range("E64").Value = 0.42
$ws.Range("F64").Value = 36
